$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set number format to Text for price cells whose new value would otherwise
# be auto-converted to a numeric type by Excel, so they remain text like the source data.
$textCells = @("D5", "D8", "D15", "D18", "D19", "D22", "D25", "D26", "D29", "D30", "D40", "D41", "D42", "D43", "D44", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update cell values per the latest crypto data refresh
$ws.Range("D2").Value = "29.906.95"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.635.00"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.90%  "
$ws.Range("D5").Value = "215.36"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").Value = "28.74"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "1.868.22"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "1.631.80"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").Value = "9.52"
$ws.Range("E15").Value = "  +7.03%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "29.915.56"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "64.70"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "240.63"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "9.88"
$ws.Range("E22").Value = "  +2.90%  "
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "157.62"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").Value = "15.53"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").Value = "0.0489"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "1.422.20"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "76.05"
$ws.Range("E40").Value = "  +9.60%  "
$ws.Range("D41").Value = "0.558"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.834"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.99"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0500"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "1.775.89"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "48.76"
$ws.Range("E49").Value = "  -9.54%  "
$ws.Range("D50").Value = "93.02"
$ws.Range("E50").Value = "  +5.50%  "
$ws.Range("D51").Value = "0.0₆0111"
$ws.Range("E51").Value = "  +10.90%  "
